# Applies: "Added labels, added dice to 4d6"
# Sets column D ("Completed") to "Yes" for a handful of rows that were
# previously blank, and moves the active selection to D21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(4, 5, 6, 11, 18, 19)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "Yes"
}

$ws.Range("D21").Select()
